$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 373, shifting the existing rows 373:392 down to 374:393
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new weekly record
$ws.Cells.Item(373, 1).Value = 3
$ws.Cells.Item(373, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(373, 3).Value = "Coquimbo"
$ws.Cells.Item(373, 4).Value = 44706
$ws.Cells.Item(373, 5).Value = 5
$ws.Cells.Item(373, 6).Value = 100112017
$ws.Cells.Item(373, 7).Value = "Apio"
$ws.Cells.Item(373, 8).Value = "Americana (o)"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 250
$ws.Cells.Item(373, 11).Value = 9000
$ws.Cells.Item(373, 12).Value = 9500
$ws.Cells.Item(373, 13).Value = 9260
$ws.Cells.Item(373, 14).Value = "`$/docena de matas"
$ws.Cells.Item(373, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(373, 16).Value = 1543
$ws.Cells.Item(373, 17).Value = 6
$ws.Cells.Item(373, 18).Value = "Hortaliza"
